$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated "FilesTab" query text (B4): the "File Type" and "Breed" coalesce
# lines were dropped from the RETURN clause (ICDC Breed script correction).
$newFilesQuery = @'
MATCH (f:file)-->(parent)
WITH DISTINCT f, parent
MATCH (f)-[*]->(c:case)<--(demo:demographic)
WHERE demo.breed IN['Greyhound'] 
OPTIONAL MATCH (s:study)<-[*]-(c)<--(diag:diagnosis)
OPTIONAL MATCH (samp:sample)-->(c)
WITH DISTINCT f, parent, c, demo, diag, s
RETURN  coalesce(f.file_name, '') AS `File Name`,
         coalesce(labels(parent)[0], '') AS `Association`,
        coalesce(f.file_description, '') AS `Description`,
        coalesce(f.file_format, '') AS `Format`,
        coalesce(f.file_size, '') AS `Size`,
        coalesce(c.case_id, '') AS `Case ID`,
         coalesce(diag.disease_term,'') AS Diagnosis , 
        coalesce(s.clinical_study_designation,'') AS `Study Code`
'@

$ws.Range("B4").Value2 = $newFilesQuery

# Row 4 is now shorter, so its autofit/explicit row height shrinks.
$ws.Rows.Item(4).RowHeight = 217.5

# Selection moved from B2 to B4 and the view scrolled so row 4 is visible.
[void]$ws.Range("B4").Select()
